$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 180, pushing the existing rows 180-263 down to 181-264.
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new record.
$ws.Range("A180").Value = 10
$ws.Range("B180").Value = "Vega Modelo de Temuco"
$ws.Range("C180").Value = "La Araucanía"
$ws.Range("D180").Value = 44845
$ws.Range("E180").Value = 9
$ws.Range("F180").Value = 100112013
$ws.Range("G180").Value = "Alcachofa"
$ws.Range("H180").Value = "Española"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 300
$ws.Range("K180").Value = 300
$ws.Range("L180").Value = 300
$ws.Range("M180").Value = 300
$ws.Range("N180").Value = "$/unidad"
$ws.Range("O180").Value = "Región Metropolitana"
$ws.Range("P180").Value = 300
$ws.Range("Q180").Value = 1
$ws.Range("R180").Value = "Hortaliza"
